$d = $word.ActiveDocument

# 1. Merge the "hitting" / " only one pin on every throw in a game" runs
#    (re-typing the italic test-one description as a single run).
$d.Content.Find.Execute("hitting only one pin on every throw in a game", $true, $false, $false, $false, $false, $true, 1, $false, "hitting only one pin on every throw in a game", 2) | Out-Null

# Re-scan paragraph indices after the text fixup above (Find/Replace does not
# change paragraph count, but we re-locate anchors by text to stay robust).
$goBackIndex = -1
$pythonDocIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Test seven:") {
        $goBackIndex = $i + 1
    }
    if ($t -match "Python Doc comments:") {
        $pythonDocIndex = $i
    }
}

# 2. Strip the stray "_GoBack" bookmark paragraph left after "Test seven:" —
#    Word drops this transient last-edit-location bookmark on a clean save,
#    leaving a bare empty paragraph behind.
if ($goBackIndex -gt 0) {
    $para = $d.Paragraphs.Item($goBackIndex)
    $para.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>") | Out-Null
}

# 3. Remove one of the three blank paragraphs between "Python Doc comments:"
#    and "Summary:" (now just two).
if ($pythonDocIndex -gt 0) {
    $blank = $d.Paragraphs.Item($pythonDocIndex + 1)
    $blank.Range.Delete() | Out-Null
}
